# Update mods data [2026-01-15 15:13:38]
# Append a new row (66) to the ModCounts sheet with the latest sample:
#   A66 = 2026/01/15 (text, same "yyyy/mm/dd" literal style as the rest of column A)
#   B66 = 逃离鸭科夫   (same game name as every other row)
#   C66 = 1144        (numeric mod count, unchanged from the previous day)
# and keep the same visual style (centered alignment) already used by the
# rest of the data rows (e.g. row 65).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 65
$newRow  = 66

# --- A66 -------------------------------------------------------------
# A plain ".Value = '2026/01/15'" assignment gets auto-parsed into a real
# date serial by the COM layer (like Excel's own "smart" entry does), which
# is NOT what the source data uses elsewhere in the column (every existing
# date cell is literal text). Writing it as a text-formula result and then
# collapsing that formula down to a static value keeps it a plain string
# without ever touching the cell's NumberFormat (which would otherwise
# leave a stray style behind).
$ws.Range("A$newRow").Formula = "=""2026/01/15"""
$ws.Range("A$newRow").Copy()
$ws.Range("A$newRow").PasteSpecial(-4163)   # xlPasteValues

# --- B66 / C66 ---------------------------------------------------------
$ws.Range("B$newRow").Value = "逃离鸭科夫"
$ws.Range("C$newRow").Value = 1144

# --- Formatting ----------------------------------------------------------
# Copy the existing row's formatting (centered alignment style) onto the
# new row so it matches the rest of the table instead of using the default
# style.
$ws.Range("A${lastRow}:C${lastRow}").Copy()
$ws.Range("A${newRow}:C${newRow}").PasteSpecial(-4122)   # xlPasteFormats

Write-Output "Added row $newRow : $($ws.Range("A$newRow").Value()) | $($ws.Range("B$newRow").Value()) | $($ws.Range("C$newRow").Value())"
